$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.494.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.349.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.352.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.10%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.881.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.346.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "63.550.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "642.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0725"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.44%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.923.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0398"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.87%  "
